# Auto update Excel log
# Appends new sensor-log rows to the PIR, Humidity and Temperature sheets.

$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $r, $c, $val) {
    $cell = $ws.Cells.Item($r, $c)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# ---------------------------------------------------------------
# PIR sheet: append rows 641-653 (Date, Timestamp, Hour, Location, Value, Status)
# ---------------------------------------------------------------
$pirWs = $wb.Worksheets.Item("PIR")
$pirData = @(
    @(641, "10:34:37", "No Motion",       "Inactive"),
    @(642, "10:34:38", "Motion Detected", "Active"),
    @(643, "10:34:45", "No Motion",       "Inactive"),
    @(644, "10:34:47", "Motion Detected", "Active"),
    @(645, "10:34:54", "No Motion",       "Inactive"),
    @(646, "10:34:55", "Motion Detected", "Active"),
    @(647, "10:35:04", "No Motion",       "Inactive"),
    @(648, "10:35:09", "No Motion",       "Inactive"),
    @(649, "10:35:14", "No Motion",       "Inactive"),
    @(650, "10:35:19", "No Motion",       "Inactive"),
    @(651, "10:35:24", "No Motion",       "Inactive"),
    @(652, "10:35:26", "No Motion",       "Inactive"),
    @(653, "10:35:34", "No Motion",       "Inactive")
)
foreach ($row in $pirData) {
    $r = $row[0]
    Set-TextCell $pirWs $r 1 "2026-02-06"
    Set-TextCell $pirWs $r 2 $row[1]
    Set-TextCell $pirWs $r 3 "10:00"
    Set-TextCell $pirWs $r 4 "Bathroom"
    Set-TextCell $pirWs $r 5 $row[2]
    Set-TextCell $pirWs $r 6 $row[3]
}

# ---------------------------------------------------------------
# Humidity sheet: append rows 460-470 (Date, Timestamp, Hour, Location, Value, Status)
# ---------------------------------------------------------------
$humWs = $wb.Worksheets.Item("Humidity")
$humData = @(
    @(460, "10:34:39", "68.3%"),
    @(461, "10:34:42", "68.1%"),
    @(462, "10:34:51", "67.7%"),
    @(463, "10:34:56", "66.6%"),
    @(464, "10:35:01", "67.5%"),
    @(465, "10:35:06", "66.7%"),
    @(466, "10:35:11", "67.8%"),
    @(467, "10:35:16", "66.7%"),
    @(468, "10:35:21", "67.5%"),
    @(469, "10:35:31", "67.4%"),
    @(470, "10:35:36", "66.5%")
)
foreach ($row in $humData) {
    $r = $row[0]
    Set-TextCell $humWs $r 1 "2026-02-06"
    Set-TextCell $humWs $r 2 $row[1]
    Set-TextCell $humWs $r 3 "10:00"
    Set-TextCell $humWs $r 4 "Bathroom"
    Set-TextCell $humWs $r 5 $row[2]
    Set-TextCell $humWs $r 6 "Active"
}

# ---------------------------------------------------------------
# Temperature sheet: append rows 459-468 (Date, Timestamp, Hour, Location, Value, Status)
# ---------------------------------------------------------------
$tempWs = $wb.Worksheets.Item("Temperature")
$tempData = @(
    @(459, "10:34:41", "28.5C"),
    @(460, "10:34:44", "28.4C"),
    @(461, "10:34:52", "28.4C"),
    @(462, "10:34:58", "28.4C"),
    @(463, "10:35:02", "28.4C"),
    @(464, "10:35:07", "28.4C"),
    @(465, "10:35:12", "28.5C"),
    @(466, "10:35:17", "28.5C"),
    @(467, "10:35:22", "28.5C"),
    @(468, "10:35:32", "28.5C")
)
foreach ($row in $tempData) {
    $r = $row[0]
    Set-TextCell $tempWs $r 1 "2026-02-06"
    Set-TextCell $tempWs $r 2 $row[1]
    Set-TextCell $tempWs $r 3 "10:00"
    Set-TextCell $tempWs $r 4 "Bathroom"
    Set-TextCell $tempWs $r 5 $row[2]
    Set-TextCell $tempWs $r 6 "Active"
}

Write-Output "Appended $($pirData.Count) PIR rows, $($humData.Count) Humidity rows, $($tempData.Count) Temperature rows."
